$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-10-21"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 10-21)"

# Update the October 2022 count (row 11) and the Total 2022 count (row 14)
$ws.Range("I11").Value = 70
$ws.Range("I14").Value = 1347
